$d = $word.ActiveDocument

$d.Content.Find.Execute("648×5=3240", $true, $false, $false, $false, $false, $true, 1, $false, "681×7=4767", 2)
$d.Content.Find.Execute("780×2=1560", $true, $false, $false, $false, $false, $true, 1, $false, "753×8=6024", 2)
$d.Content.Find.Execute("279×3=837", $true, $false, $false, $false, $false, $true, 1, $false, "320×4=1280", 2)
$d.Content.Find.Execute("661×5=3305", $true, $false, $false, $false, $false, $true, 1, $false, "704×5=3520", 2)
$d.Content.Find.Execute("815×6=4890", $true, $false, $false, $false, $false, $true, 1, $false, "443×9=3987", 2)
$d.Content.Find.Execute("382×6=2292", $true, $false, $false, $false, $false, $true, 1, $false, "938×5=4690", 2)
$d.Content.Find.Execute("142×8=1136", $true, $false, $false, $false, $false, $true, 1, $false, "845×7=5915", 2)
$d.Content.Find.Execute("390×5=1950", $true, $false, $false, $false, $false, $true, 1, $false, "998×2=1996", 2)
$d.Content.Find.Execute("326×3=978", $true, $false, $false, $false, $false, $true, 1, $false, "381×8=3048", 2)
$d.Content.Find.Execute("545×9=4905", $true, $false, $false, $false, $false, $true, 1, $false, "624×6=3744", 2)
$d.Content.Find.Execute("735×3=2205", $true, $false, $false, $false, $false, $true, 1, $false, "637×3=1911", 2)
$d.Content.Find.Execute("816×6=4896", $true, $false, $false, $false, $false, $true, 1, $false, "128×7=896", 2)
$d.Content.Find.Execute("277×8=2216", $true, $false, $false, $false, $false, $true, 1, $false, "551×6=3306", 2)
$d.Content.Find.Execute("972×4=3888", $true, $false, $false, $false, $false, $true, 1, $false, "610×3=1830", 2)
$d.Content.Find.Execute("296×9=2664", $true, $false, $false, $false, $false, $true, 1, $false, "845×5=4225", 2)
$d.Content.Find.Execute("354×2=708", $true, $false, $false, $false, $false, $true, 1, $false, "379×6=2274", 2)
$d.Content.Find.Execute("764×2=1528", $true, $false, $false, $false, $false, $true, 1, $false, "937×9=8433", 2)
$d.Content.Find.Execute("120×6=720", $true, $false, $false, $false, $false, $true, 1, $false, "961×2=1922", 2)
$d.Content.Find.Execute("171×9=1539", $true, $false, $false, $false, $false, $true, 1, $false, "901×5=4505", 2)
$d.Content.Find.Execute("607×3=1821", $true, $false, $false, $false, $false, $true, 1, $false, "736×8=5888", 2)
$d.Content.Find.Execute("962×8=7696", $true, $false, $false, $false, $false, $true, 1, $false, "846×3=2538", 2)
$d.Content.Find.Execute("842×3=2526", $true, $false, $false, $false, $false, $true, 1, $false, "973×5=4865", 2)
$d.Content.Find.Execute("736×9=6624", $true, $false, $false, $false, $false, $true, 1, $false, "997×7=6979", 2)
$d.Content.Find.Execute("772×5=3860", $true, $false, $false, $false, $false, $true, 1, $false, "909×2=1818", 2)
$d.Content.Find.Execute("624×9=5616", $true, $false, $false, $false, $false, $true, 1, $false, "296×7=2072", 2)
